$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (within the used data range only, to
# avoid materializing the entire 1,048,576-row column), shifting the
# existing D:K figures one column to the right (now E:L).
$ws.Range("D7:D102").Insert(-4161)

# Copy formatting (number formats/styles) from the (now shifted) column E
# into the newly inserted column D so the new column matches the original
# per-row styling (date format rows 7/38/80, numeric format elsewhere).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column D with the latest reporting period's figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 380400
$ws.Range("D9").Value = 119600
$ws.Range("D10").Value = 260800
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 96300
$ws.Range("D17").Value = 263200
$ws.Range("D18").Value = 117200
$ws.Range("D20").Value = 13600
$ws.Range("D21").Value = 227000
$ws.Range("D22").Value = 63500
$ws.Range("D23").Value = 67200
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 67200
$ws.Range("D27").Value = 36700
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -13600
$ws.Range("D33").Value = 36700
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 36700
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 113700
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 16800
$ws.Range("D44").Value = 11000
$ws.Range("D45").Value = 29300
$ws.Range("D46").Value = 170800
$ws.Range("D47").Value = 148900
$ws.Range("D48").Value = 2608700
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 122500
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 3050800
$ws.Range("D57").Value = 8600
$ws.Range("D58").Value = 183500
$ws.Range("D59").Value = 32600
$ws.Range("D60").Value = 224700
$ws.Range("D61").Value = 1464300
$ws.Range("D62").Value = 4700
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1693700
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 38700
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1357100
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 36700
$ws.Range("D83").Value = 96300
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 140800
$ws.Range("D91").Value = -143000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -112600
$ws.Range("D96").Value = -49100
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -80500
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -52400

